# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") on Sheet1 is recalculated from source data; this script
# writes the recomputed K values for data rows 2-50 (column G, 7th column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value (column G), for rows 2 through 50.
$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 3
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 3
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 4
    38 = 1
    39 = 0
    40 = 3
    41 = 1
    42 = 1
    43 = 2
    44 = 0
    45 = 0
    46 = 3
    47 = 1
    48 = 1
    49 = 0
    50 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
